$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.257.53'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").Value = '3.497.51'
$ws.Range("E3").Value = '  +0.91%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.11'
$ws.Range("D5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.95'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.92%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("E9").Value = '  +2.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.20'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.73%  '

$ws.Range("E11").Value = '  +1.87%  '

$ws.Range("D12").Value = '4.095.53'
$ws.Range("E12").Value = '  +0.84%  '

$ws.Range("E13").Value = '  +1.36%  '

$ws.Range("E14").Value = '  +2.86%  '

$ws.Range("D15").Value = '3.498.40'
$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.07%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '64.288.94'
$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.93'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.78'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.73'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '386.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").Value = '3.637.38'
$ws.Range("E23").Value = '  +0.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.38%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  +1.58%  '

$ws.Range("E27").Value = '  +3.18%  '

$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.41%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.44'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.29'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.23'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.47'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -6.24%  '

$ws.Range("D33").Value = '3.518.05'
$ws.Range("E33").Value = '  +1.38%  '

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.149'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.53%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.55'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.20'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.55'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.88'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.42'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0780'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.804'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.85'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.29%  '

$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.77'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.41'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.82%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.65'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.13%  '

$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.17'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.28%  '

$ws.Range("D49").Value = '2.480.23'
$ws.Range("E49").Value = '  +2.68%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.76'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.897'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.06%  '
